$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell text updates (crypto price / volume refresh + two rank swaps).
# Every new value is written with a leading apostrophe (quote-prefix) so
# Excel keeps numeric-looking strings (e.g. "507.06") as plain text,
# matching the source inline-string cells, instead of silently coercing
# them to numbers. The cell style is then reset to "Normal" so no stray
# quote-prefix / @ text-format style is left behind afterwards.
$ws.Range('D2').Value = "'" + '56.581.51'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -3.62%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '2.375.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -4.60%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.24%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '507.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -4.98%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '129.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -2.82%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.30%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.556'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -2.03%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '2.397.57'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -4.07%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.0967'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -2.36%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -1.63%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '0.324'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -1.44%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '4.68'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -9.97%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '2.800.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -4.51%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '56.410.94'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -3.82%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -2.54%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -2.78%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '2.394.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -4.09%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -3.03%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '313.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -2.02%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -4.01%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '6.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +0.75%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -0.01%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '65.94'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +0.02%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +0.31%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '2.494.42'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -4.25%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '0.377'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -7.48%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -4.48%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '7.24'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -2.76%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '174.59'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +0.92%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = "'" + 'PEPE'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'" + 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'" + '0.0₃0718'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -4.83%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'" + 'PancakeSwap'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'" + 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'" + '1.67'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -3.59%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '6.16'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -1.50%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '1.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -5.99%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +0.00%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '0.995'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -0.20%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '17.79'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -1.57%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -1.03%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -4.62%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '35.83'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -1.19%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -4.62%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '0.788'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -2.35%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '133.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +1.85%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '3.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -2.63%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '4.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -5.16%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '256.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -6.53%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -3.18%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -3.29%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B50').Value = "'" + 'EnergySwap'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'" + '16.84'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -4.05%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'" + 'VeChain'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'" + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'" + '0.0208'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -4.50%  '
$ws.Range('E51').Style = 'Normal'
